$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2647968.5
$ws.Cells.Item(17, 9).Value = 1263
$ws.Cells.Item(17, 11).Value = 3789
$ws.Cells.Item(17, 13).Value = -3621

$ws.Cells.Item(19, 8).Value = 1480.6428
$ws.Cells.Item(19, 9).Value = 979.4
$ws.Cells.Item(19, 10).Value = 1759.1111
$ws.Cells.Item(19, 11).Value = 979.4
$ws.Cells.Item(19, 12).Value = 1759.1111
$ws.Cells.Item(19, 13).Value = -804.4
$ws.Cells.Item(19, 14).Value = -2109.1111

$ws.Cells.Item(62, 8).Value = 6906500
$ws.Cells.Item(62, 9).Value = 14293464
$ws.Cells.Item(62, 10).Value = 12000
$ws.Cells.Item(62, 11).Value = 14293464
$ws.Cells.Item(62, 12).Value = 12000
$ws.Cells.Item(62, 13).Value = -14292840
$ws.Cells.Item(62, 14).Value = -13248

$ws.Cells.Item(65, 8).Value = 6906500
$ws.Cells.Item(65, 9).Value = 14293464
$ws.Cells.Item(65, 10).Value = 12000
$ws.Cells.Item(65, 11).Value = 71467320
$ws.Cells.Item(65, 12).Value = 60000
$ws.Cells.Item(65, 13).Value = -71464200
$ws.Cells.Item(65, 14).Value = -66240

$ws.Cells.Item(112, 8).Value = 68194.92999999999
$ws.Cells.Item(112, 10).Value = 101981
$ws.Cells.Item(112, 12).Value = 305943
$ws.Cells.Item(112, 14).Value = -308159

$ws.Cells.Item(141, 8).Value = 1935.2
$ws.Cells.Item(141, 9).Value = 1773.8948
$ws.Cells.Item(141, 11).Value = 5321.6844
$ws.Cells.Item(141, 13).Value = -141.6844000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 28574408
$ws.Cells.Item(2, 9).Value = 33335940
$ws.Cells.Item(2, 11).Value = 33335940
$ws.Cells.Item(2, 13).Value = -33335827

$ws.Cells.Item(88, 8).Value = 2206.8333
$ws.Cells.Item(88, 9).Value = 2137.3333
$ws.Cells.Item(88, 10).Value = 2276.3333
$ws.Cells.Item(88, 11).Value = 2137.3333
$ws.Cells.Item(88, 12).Value = 2276.3333
$ws.Cells.Item(88, 14).Value = -3088.3333
$ws.Cells.Item(88, 13).Value = -1731.3333

$ws.Cells.Item(91, 8).Value = 2206.8333
$ws.Cells.Item(91, 9).Value = 2137.3333
$ws.Cells.Item(91, 10).Value = 2276.3333
$ws.Cells.Item(91, 11).Value = 2137.3333
$ws.Cells.Item(91, 12).Value = 2276.3333
$ws.Cells.Item(91, 14).Value = -5084.3333
$ws.Cells.Item(91, 13).Value = -733.3332999999998

$ws.Cells.Item(102, 8).Value = 5307.067
$ws.Cells.Item(102, 9).Value = 4420.9
$ws.Cells.Item(102, 10).Value = 7079.4
$ws.Cells.Item(102, 11).Value = 4420.9
$ws.Cells.Item(102, 12).Value = 7079.4
$ws.Cells.Item(102, 13).Value = -2798.9
$ws.Cells.Item(102, 14).Value = -10323.4

$ws.Cells.Item(116, 8).Value = 28574408
$ws.Cells.Item(116, 9).Value = 33335940
$ws.Cells.Item(116, 11).Value = 33335940
$ws.Cells.Item(116, 13).Value = -33333646

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 28574408
$ws.Cells.Item(3, 9).Value = 33335940
$ws.Cells.Item(3, 11).Value = 33335940
$ws.Cells.Item(3, 13).Value = -33335826

$ws.Cells.Item(105, 8).Value = 1500.2354
$ws.Cells.Item(105, 9).Value = 1366.9333
$ws.Cells.Item(105, 11).Value = 1366.9333
$ws.Cells.Item(105, 13).Value = 380.0667000000001

$ws.Cells.Item(134, 8).Value = 1876.5405
$ws.Cells.Item(134, 9).Value = 1887
$ws.Cells.Item(134, 10).Value = 1500
$ws.Cells.Item(134, 11).Value = 5661
$ws.Cells.Item(134, 12).Value = 4500
$ws.Cells.Item(134, 13).Value = -3126
$ws.Cells.Item(134, 14).Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 8588.134
$ws.Cells.Item(99, 9).Value = 8164.6665
$ws.Cells.Item(99, 11).Value = 8164.6665
$ws.Cells.Item(99, 13).Value = -6666.6665

$ws.Cells.Item(105, 8).Value = 1110.3334
$ws.Cells.Item(105, 9).Value = 1166.5
$ws.Cells.Item(105, 10).Value = 998
$ws.Cells.Item(105, 11).Value = 1166.5
$ws.Cells.Item(105, 12).Value = 998
$ws.Cells.Item(105, 13).Value = 580.5
$ws.Cells.Item(105, 14).Value = -4492

$ws.Cells.Item(126, 8).Value = 8588.134
$ws.Cells.Item(126, 9).Value = 8164.6665
$ws.Cells.Item(126, 11).Value = 24493.9995
$ws.Cells.Item(126, 13).Value = -22023.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 305.5
$ws.Cells.Item(40, 10).Value = 574.2857
$ws.Cells.Item(40, 12).Value = 2297.1428
$ws.Cells.Item(40, 14).Value = -2435.1428

$ws.Cells.Item(55, 8).Value = 150
$ws.Cells.Item(55, 9).Value = 150
$ws.Cells.Item(55, 11).Value = 450
$ws.Cells.Item(55, 13).Value = -273

$ws.Cells.Item(129, 8).Value = 733.4545000000001
$ws.Cells.Item(129, 9).Value = 671
$ws.Cells.Item(129, 10).Value = 900
$ws.Cells.Item(129, 11).Value = 2013
$ws.Cells.Item(129, 12).Value = 2700
$ws.Cells.Item(129, 13).Value = 2987
$ws.Cells.Item(129, 14).Value = -12700

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 11995.19
$ws.Cells.Item(43, 9).Value = 9272.166999999999
$ws.Cells.Item(43, 10).Value = 28333.334
$ws.Cells.Item(43, 11).Value = 9272.166999999999
$ws.Cells.Item(43, 12).Value = 28333.334
$ws.Cells.Item(43, 13).Value = -9121.166999999999
$ws.Cells.Item(43, 14).Value = -28635.334

$ws.Cells.Item(122, 8).Value = 2508.2856
$ws.Cells.Item(122, 9).Value = 1569.6666
$ws.Cells.Item(122, 10).Value = 3212.25
$ws.Cells.Item(122, 11).Value = 4708.9998
$ws.Cells.Item(122, 12).Value = 9636.75
$ws.Cells.Item(122, 13).Value = -2258.9998
$ws.Cells.Item(122, 14).Value = -14536.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 14207.223
$ws.Cells.Item(7, 9).Value = 17647.666
$ws.Cells.Item(7, 10).Value = 7326.3335
$ws.Cells.Item(7, 11).Value = 17647.666
$ws.Cells.Item(7, 12).Value = 7326.3335
$ws.Cells.Item(7, 13).Value = -17535.666
$ws.Cells.Item(7, 14).Value = -7550.3335

$ws.Cells.Item(16, 8).Value = 10589.8
$ws.Cells.Item(16, 9).Value = 4950
$ws.Cells.Item(16, 11).Value = 4950
$ws.Cells.Item(16, 13).Value = -4780

$ws.Cells.Item(42, 8).Value = 12919.23
$ws.Cells.Item(42, 9).Value = 11995
$ws.Cells.Item(42, 10).Value = 13330
$ws.Cells.Item(42, 11).Value = 11995
$ws.Cells.Item(42, 12).Value = 13330
$ws.Cells.Item(42, 14).Value = -14456
$ws.Cells.Item(42, 13).Value = -11432

$ws.Cells.Item(49, 8).Value = 12919.23
$ws.Cells.Item(49, 9).Value = 11995
$ws.Cells.Item(49, 10).Value = 13330
$ws.Cells.Item(49, 11).Value = 11995
$ws.Cells.Item(49, 12).Value = 13330
$ws.Cells.Item(49, 14).Value = -13624
$ws.Cells.Item(49, 13).Value = -11848

$ws.Cells.Item(68, 8).Value = 4333.2
$ws.Cells.Item(68, 10).Value = 6333
$ws.Cells.Item(68, 12).Value = 6333
$ws.Cells.Item(68, 14).Value = -7831

$ws.Cells.Item(71, 8).Value = 4333.2
$ws.Cells.Item(71, 10).Value = 6333
$ws.Cells.Item(71, 12).Value = 31665
$ws.Cells.Item(71, 14).Value = -39153

$ws.Cells.Item(100, 8).Value = 99.5
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 14207.223
$ws.Cells.Item(126, 9).Value = 17647.666
$ws.Cells.Item(126, 10).Value = 7326.3335
$ws.Cells.Item(126, 11).Value = 52942.99800000001
$ws.Cells.Item(126, 12).Value = 21979.0005
$ws.Cells.Item(126, 13).Value = -50472.99800000001
$ws.Cells.Item(126, 14).Value = -26919.0005

$ws.Cells.Item(131, 8).Value = 28666.666
$ws.Cells.Item(131, 9).Value = 44000
$ws.Cells.Item(131, 10).Value = 26750
$ws.Cells.Item(131, 11).Value = 44000
$ws.Cells.Item(131, 12).Value = 26750
$ws.Cells.Item(131, 13).Value = -38960
$ws.Cells.Item(131, 14).Value = -36830

$ws.Cells.Item(132, 8).Value = 3584.2
$ws.Cells.Item(132, 9).Value = 3200.3333
$ws.Cells.Item(132, 10).Value = 4160
$ws.Cells.Item(132, 11).Value = 9600.999899999999
$ws.Cells.Item(132, 12).Value = 12480
$ws.Cells.Item(132, 13).Value = -7070.999899999999
$ws.Cells.Item(132, 14).Value = -17540

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 46564.715
$ws.Cells.Item(41, 9).Value = 30192.334
$ws.Cells.Item(41, 11).Value = 30192.334
$ws.Cells.Item(41, 13).Value = -29802.334

$ws.Cells.Item(125, 8).Value = 200022780
$ws.Cells.Item(125, 10).Value = 200022780
$ws.Cells.Item(125, 12).Value = 200022780
$ws.Cells.Item(125, 14).Value = -200032620
